$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = 100
$ws.Range("B7").Value = 350
$ws.Range("B14").Value = -200
$ws.Range("B15").Value = -250
$ws.Range("B16").Value = 1.6
$ws.Range("B21").Value = 45

$ws.Range("B16").Select()
